# Commit: change "update" to "all" in git_push
#
# Adds a new row to the "Code Executions" table (header row:
# "Project" | "URL" | "Comment") recording the commit that pushes
# newly created files:
#   Logging Notebook |
#   https://github.com/EmaadKhwaja/Word-Logging-Script/commit/5b757a35d1aaf5adee17e185e4405333618aade3 |
#   push newly created files

$d = $word.ActiveDocument

# Locate the "Code Executions" table robustly by its header row
# ("Project" / "URL" / "Comment") rather than a hard-coded index.
$target = $null
foreach ($tbl in $d.Tables) {
    if ($tbl.Columns.Count -eq 3 -and $tbl.Rows.Count -ge 1) {
        $col1 = $tbl.Cell(1, 1).Range.Text
        $col2 = $tbl.Cell(1, 2).Range.Text
        if ($col1 -like "Project*" -and $col2 -like "URL*") {
            $target = $tbl
        }
    }
}

if ($target -ne $null) {
    $newRow = $target.Rows.Add()
    $idx = $newRow.Index

    $target.Cell($idx, 1).Range.Text = "Logging Notebook"
    $target.Cell($idx, 2).Range.Text = "https://github.com/EmaadKhwaja/Word-Logging-Script/commit/5b757a35d1aaf5adee17e185e4405333618aade3"
    $target.Cell($idx, 3).Range.Text = "push newly created files"
}
